# Automatische test-sync: 2025-08-19 21:01:50
# Adds the new incoming mail-log entry to the "Logs" sheet and bumps the
# matching tally on the "Dashboard" sheet.

$wb = $excel.ActiveWorkbook

$logs = $wb.Worksheets.Item("Logs")
$dashboard = $wb.Worksheets.Item("Dashboard")

$newRow = 21

$logs.Cells.Item($newRow, 1).Value = "Vraag over product"
$logs.Cells.Item($newRow, 2).Value = "documentatie@testbedrijf123.nl"
$logs.Cells.Item($newRow, 4).Value = "Intern verzoek / Actie voor medewerker"
$logs.Cells.Item($newRow, 6).Value = "2025-08-19 21:01:05"
$logs.Cells.Item($newRow, 7).Value = "Nee"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# Extend the conditional-formatting ranges so they keep covering the whole
# data range (D/G/H/I/J2:20 -> 2:21), same as Excel does when you drag the
# formatted range's border down over a newly filled-in row.
$logs.Range("D2:D20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D21"))
$logs.Range("G2:G20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G21"))
$logs.Range("H2:H20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H21"))
$logs.Range("I2:I20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I21"))
$logs.Range("J2:J20").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J21"))

# Bump the dashboard tally for this category (19 -> 20).
$dashboard.Range("B2").Value = 20
